$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the row for "Полубоярцев Максим" (was row 10), who has no numeric
# grades (only "Н"/"H" placeholders). This shifts rows 11-13 up to 10-12.
$ws.Rows(10).Delete()

# Re-apply the sort so the sortState range shrinks to match the new data
# extent (A2:E12) instead of the stale A2:E13.
$sortObj = $ws.Sort
$sortObj.SortFields.Clear()
$sortObj.SortFields.Add($ws.Range("A2:A12")) | Out-Null
$sortObj.SetRange($ws.Range("A2:E12"))
$sortObj.Header = -4142
$sortObj.Apply()

# New column P header for task №15
$ws.Range("P1").Value = "№15"

# Fill in grades for task №13 (column N) for every student
$ws.Range("N2").Value = 5
$ws.Range("N3").Value = 5
$ws.Range("N4").Value = 0
$ws.Range("N5").Value = 5
$ws.Range("N6").Value = 0
$ws.Range("N7").Value = 5
$ws.Range("N8").Value = 0
$ws.Range("N9").Value = 5
$ws.Range("N10").Value = 0
$ws.Range("N11").Value = 5
$ws.Range("N12").Value = 0

# New grade for task №14 (column O) for Домнин Александр
$ws.Range("O7").Value = 5

# Misc grade corrections for earlier tasks
$ws.Range("F4").Value = "4(late)"
$ws.Range("L5").Value = "4(late)"
$ws.Range("M5").Value = "4(late)"
$ws.Range("M6").Value = "4(late)"
$ws.Range("L8").Value = "4(late)"
$ws.Range("M8").Value = "4(late)"
$ws.Range("J10").Value = "4(late)"
$ws.Range("K10").Value = "4(late)"
$ws.Range("L10").Value = "4(late)"
$ws.Range("M10").Value = "4(late)"

# Highlight Хващевский Дмитрий's name in red, like Владимир Шекуров
$ws.Range("A12").Font.Color = 255

# Page setup: portrait A4
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1

# Update the active cell/selection to P4
$ws.Range("P4").Select() | Out-Null
